$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.875.11'
$ws.Range("E2").Value = '  +1.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.271.56'
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.28'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.82'
$ws.Range("E6").Value = '  +1.18%  '
$ws.Range("E7").Value = '  +2.02%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.60'
$ws.Range("E10").Value = '  +1.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.54'
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0797'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.624.10'
$ws.Range("E15").Value = '  +0.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.31'
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.286.60'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.778'
$ws.Range("E18").Value = '  +3.71%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.786.13'
$ws.Range("E19").Value = '  +1.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.66'
$ws.Range("E20").Value = '  +4.09%  '
$ws.Range("E21").Value = '  +0.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.95'
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.20'
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '243.50'
$ws.Range("E24").Value = '  +1.60%  '
$ws.Range("E25").Value = '  +0.61%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("B27").Value = 'ImmutableX'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.93'
$ws.Range("E27").Value = '  +3.57%  '
$ws.Range("E28").Value = '  +1.57%  '
$ws.Range("E29").Value = '  -1.10%  '
$ws.Range("E30").Value = '  -1.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.34'
$ws.Range("E31").Value = '  +4.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '160.92'
$ws.Range("E32").Value = '  +2.07%  '
$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").Value = '  +1.13%  '
$ws.Range("E36").Value = '  -0.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.10'
$ws.Range("E37").Value = '  +3.27%  '
$ws.Range("E38").Value = '  +2.39%  '
$ws.Range("E39").Value = '  -0.18%  '
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("E41").Value = '  +1.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.94'
$ws.Range("E42").Value = '  -1.22%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.008.66'
$ws.Range("E43").Value = '  -2.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.55'
$ws.Range("E44").Value = '  -3.14%  '
$ws.Range("E45").Value = '  +1.84%  '
$ws.Range("E46").Value = '  +1.91%  '
$ws.Range("E47").Value = '  +2.15%  '
$ws.Range("E48").Value = '  -1.83%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '52.90'
$ws.Range("E49").Value = '  +3.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.52'
$ws.Range("E50").Value = '  +0.51%  '
$ws.Range("E51").Value = '  +1.19%  '
